# Avantis Mapping / SPARQL_superclass.xlsx
# Commit: "work on prep for removal" - Worked on class-less assets and classes to remove

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the plain-text class labels in column A into full TWONTO URIs,
# and fix up the corresponding Avantis class in column B for the rows
# that were mis-mapped (air duct segment <-> instrument gauge or display).

$ws.Range("A3").Value = "http://www.toronto.ca/TWONTO#air_duct_segment"
$ws.Range("B3").Value = "Silencer"

$ws.Range("A4").Value = "http://www.toronto.ca/TWONTO#cable_segment"
$ws.Range("B4").Value = "Electrical Power Line"

$ws.Range("A5").Value = "http://www.toronto.ca/TWONTO#instrument_gauge_or_display"
$ws.Range("B5").Value = "Pressure Indicator"

# Update the current selection to match the author's last interaction.
$ws.Range("A3:C5").Select()
